# Insert a new row into the language sheet for the MAGMA chamber title,
# pushing the existing rows (9 onward) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9; existing row 9 (and below) shift to row 10 (and below).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row with the new key/value pair.
$ws.Cells.Item(9, 1).Value = "magma_title"
$ws.Cells.Item(9, 2).Value = "MAGMA"

# Move the active selection to reflect where the new row was added.
$ws.Range("B9").Select()
